$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 528.6667
$ws.Range("H40").Value = 1085
$ws.Range("I40").Value = 760
$ws.Range("J40").Value = 1280
$ws.Range("K40").Value = 760
$ws.Range("L40").Value = 1280
$ws.Range("M40").Value = -585
$ws.Range("N40").Value = -1630
$ws.Range("H43").Value = 3395.85
$ws.Range("I43").Value = 853.3333
$ws.Range("J43").Value = 4485.5
$ws.Range("K43").Value = 853.3333
$ws.Range("L43").Value = 4485.5
$ws.Range("M43").Value = -784.3333
$ws.Range("N43").Value = -4623.5
$ws.Range("H129").Value = 250949.66
$ws.Range("J129").Value = 271248.8
$ws.Range("L129").Value = 813746.3999999999
$ws.Range("N129").Value = -823746.3999999999
$ws.Range("H131").Value = 3121.5217
$ws.Range("I131").Value = 1448.75
$ws.Range("J131").Value = 3473.6843
$ws.Range("K131").Value = 4346.25
$ws.Range("L131").Value = 10421.0529
$ws.Range("M131").Value = 693.75
$ws.Range("N131").Value = -20501.0529
$ws.Range("H138").Value = 2872.6355
$ws.Range("I138").Value = 2032.0834
$ws.Range("J138").Value = 3152.8193
$ws.Range("K138").Value = 6096.2502
$ws.Range("L138").Value = 9458.457900000001
$ws.Range("M138").Value = -956.2502000000004
$ws.Range("N138").Value = -19738.4579
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 201
$ws.Range("I4").Value = 201
$ws.Range("K4").Value = 201
$ws.Range("M4").Value = -85
$ws.Range("H5").Value = 107.75
$ws.Range("I5").Value = 107.75
$ws.Range("K5").Value = 107.75
$ws.Range("M5").Value = 4.25
$ws.Range("H32").Value = 15092.062
$ws.Range("I32").Value = 11129.375
$ws.Range("K32").Value = 11129.375
$ws.Range("M32").Value = -10842.375
$ws.Range("H63").Value = 15626500
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 15626500
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H80").Value = 41604.09
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 45464.5
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 45464.5
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -47460.5
$ws.Range("H83").Value = 41604.09
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 45464.5
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 136393.5
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -146377.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 107.75
$ws.Range("I4").Value = 107.75
$ws.Range("K4").Value = 107.75
$ws.Range("M4").Value = 7.25
$ws.Range("H99").Value = 1411
$ws.Range("I99").Value = 1501.25
$ws.Range("K99").Value = 1501.25
$ws.Range("M99").Value = -3.25
$ws.Range("H118").Value = 33333
$ws.Range("J118").Value = 33333
$ws.Range("L118").Value = 33333
$ws.Range("N118").Value = -36647
$ws.Range("H134").Value = 4199.5127
$ws.Range("I134").Value = 4182.0347
$ws.Range("J134").Value = 4250.2
$ws.Range("K134").Value = 12546.1041
$ws.Range("L134").Value = 12750.6
$ws.Range("M134").Value = -10011.1041
$ws.Range("N134").Value = -17820.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 62.333332
$ws.Range("I7").Value = 92.75
$ws.Range("K7").Value = 92.75
$ws.Range("M7").Value = 20.25
$ws.Range("H15").Value = 554
$ws.Range("I15").Value = 554
$ws.Range("K15").Value = 554
$ws.Range("M15").Value = -384
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2727407.2
$ws.Range("I4").Value = 122.375
$ws.Range("J4").Value = 10000167
$ws.Range("K4").Value = 367.125
$ws.Range("L4").Value = 30000501
$ws.Range("M4").Value = -255.125
$ws.Range("N4").Value = -30000725
$ws.Range("H86").Value = 360.66666
$ws.Range("I86").Value = 291
$ws.Range("K86").Value = 873
$ws.Range("M86").Value = 313
$ws.Range("H89").Value = 360.66666
$ws.Range("I89").Value = 291
$ws.Range("K89").Value = 2619
$ws.Range("M89").Value = 3309
$ws.Range("H131").Value = 731.62
$ws.Range("I131").Value = 407.77777
$ws.Range("J131").Value = 763.6484
$ws.Range("K131").Value = 1223.33331
$ws.Range("L131").Value = 2290.9452
$ws.Range("M131").Value = 3816.66669
$ws.Range("N131").Value = -12370.9452
$ws.Range("H140").Value = 2760.36
$ws.Range("I140").Value = 1559.3529
$ws.Range("J140").Value = 5312.5
$ws.Range("K140").Value = 4678.0587
$ws.Range("L140").Value = 15937.5
$ws.Range("M140").Value = 501.9412999999995
$ws.Range("N140").Value = -26297.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4242.1665
$ws.Range("I80").Value = 4050
$ws.Range("K80").Value = 4050
$ws.Range("M80").Value = -3052
$ws.Range("H83").Value = 4242.1665
$ws.Range("I83").Value = 4050
$ws.Range("K83").Value = 20250
$ws.Range("M83").Value = -15258
$ws.Range("H107").Value = 730
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 782
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 782
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -4622
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5600.5
$ws.Range("I22").Value = 5600.5
$ws.Range("K22").Value = 5600.5
$ws.Range("M22").Value = -5305.5
$ws.Range("H27").Value = 5600.5
$ws.Range("I27").Value = 5600.5
$ws.Range("K27").Value = 5600.5
$ws.Range("M27").Value = -5493.5
$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 1200
$ws.Range("L46").Value = 1200
$ws.Range("N46").Value = -1576
$ws.Range("H82").Value = 3200.7273
$ws.Range("I82").Value = 2200.3333
$ws.Range("J82").Value = 4401.2
$ws.Range("K82").Value = 2200.3333
$ws.Range("L82").Value = 4401.2
$ws.Range("M82").Value = -1839.3333
$ws.Range("N82").Value = -5123.2
$ws.Range("H85").Value = 3200.7273
$ws.Range("I85").Value = 2200.3333
$ws.Range("J85").Value = 4401.2
$ws.Range("K85").Value = 2200.3333
$ws.Range("L85").Value = 4401.2
$ws.Range("M85").Value = -952.3332999999998
$ws.Range("N85").Value = -6897.2
$ws.Range("H100").Value = 2016.1578
$ws.Range("I100").Value = 1340.3
$ws.Range("J100").Value = 2767.111
$ws.Range("K100").Value = 1340.3
$ws.Range("L100").Value = 2767.111
$ws.Range("M100").Value = -799.3
$ws.Range("N100").Value = -3849.111
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 60478.715
$ws.Range("I113").Value = 76654.73
$ws.Range("J113").Value = 1166.6666
$ws.Range("K113").Value = 229964.19
$ws.Range("L113").Value = 3499.9998
$ws.Range("M113").Value = -227794.19
$ws.Range("N113").Value = -7839.9998
$ws.Range("H132").Value = 12821815
$ws.Range("I132").Value = 16129972
$ws.Range("K132").Value = 48389916
$ws.Range("M132").Value = -48387386
